$wb = $excel.ActiveWorkbook

# Sheet1 (Hoja1): update value, selection, and remove tab-selected state
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("B2").Value = 122232
$ws1.Range("B9").Select()

# Sheet2 (Hoja2): update values, selection; becomes the active/selected tab
$ws2 = $wb.Worksheets.Item("Hoja2")
$ws2.Range("B2").Value = 232323
$ws2.Range("B3").Value = 10774750
$ws2.Activate()
$ws2.Range("B10").Select()
